$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at row 100, pushing existing rows 100-115 down to 102-117
$ws.Rows.Item(100).Resize(2).Insert()

# Row 100: new "Primera" quality entry for date 2022-11-11 (serial 44876)
$ws.Cells.Item(100, 1).Value = 7
$ws.Cells.Item(100, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(100, 3).Value = "Ñuble"
$ws.Cells.Item(100, 4).Value = 44876
$ws.Cells.Item(100, 5).Value = 16
$ws.Cells.Item(100, 6).Value = 100112040
$ws.Cells.Item(100, 7).Value = "Cilantro"
$ws.Cells.Item(100, 8).Value = "Sin especificar"
$ws.Cells.Item(100, 9).Value = "Primera"
$ws.Cells.Item(100, 10).Value = 400
$ws.Cells.Item(100, 11).Value = 600
$ws.Cells.Item(100, 12).Value = 700
$ws.Cells.Item(100, 13).Value = 650
$ws.Cells.Item(100, 14).Value = "`$/atado 0,5 a 1 kilo"
$ws.Cells.Item(100, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(100, 16).Value = 650
$ws.Cells.Item(100, 17).Value = 1
$ws.Cells.Item(100, 18).Value = "Hortaliza"

# Row 101: new "Segunda" quality entry for date 2022-11-11 (serial 44876)
$ws.Cells.Item(101, 1).Value = 7
$ws.Cells.Item(101, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(101, 3).Value = "Ñuble"
$ws.Cells.Item(101, 4).Value = 44876
$ws.Cells.Item(101, 5).Value = 16
$ws.Cells.Item(101, 6).Value = 100112040
$ws.Cells.Item(101, 7).Value = "Cilantro"
$ws.Cells.Item(101, 8).Value = "Sin especificar"
$ws.Cells.Item(101, 9).Value = "Segunda"
$ws.Cells.Item(101, 10).Value = 300
$ws.Cells.Item(101, 11).Value = 500
$ws.Cells.Item(101, 12).Value = 500
$ws.Cells.Item(101, 13).Value = 500
$ws.Cells.Item(101, 14).Value = "`$/atado 0,5 a 1 kilo"
$ws.Cells.Item(101, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(101, 16).Value = 500
$ws.Cells.Item(101, 17).Value = 1
$ws.Cells.Item(101, 18).Value = "Hortaliza"
